$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Metazoa")
$ws.Range("B2").Value = -0.2870643238427826
$ws.Range("C2").Value = -1.213958938283074
$ws.Range("D2").Value = 0.6571227527281447
$ws.Range("E2").Value = 0.4841473864610111
$ws.Range("G2").Value = 0.01518785678901663
$ws.Range("H2").Value = 0.0000000004336090312679973
$ws.Range("I2").Value = 0.1114465384825023
$ws.Range("B3").Value = 0.06246785514814768
$ws.Range("C3").Value = -0.3240844722207426
$ws.Range("D3").Value = 0.5243877783215349
$ws.Range("E3").Value = 0.7527849185946873
$ws.Range("G3").Value = 0.3181110463215805
$ws.Range("H3").Value = 0.03795643452122212
$ws.Range("I3").Value = 0.5479150678513206
$ws.Range("B4").Value = 0.6060071755630394
$ws.Range("C4").Value = 0.004055495560704211
$ws.Range("D4").Value = 1.192054586825716
$ws.Range("B5").Value = 0.5186528085313554
$ws.Range("C5").Value = 0.00000560036937795076
$ws.Range("D5").Value = 1.49585596558568

$ws = $wb.Worksheets.Item("Fungi")
$ws.Range("B2").Value = -2.069507545374689
$ws.Range("C2").Value = -3.295465335415409
$ws.Range("D2").Value = -0.8359086681387676
$ws.Range("E2").Value = 0.006855184233076184
$ws.Range("G2").Value = 0.01507050511751816
$ws.Range("H2").Value = 0.0000000002650000550675009
$ws.Range("I2").Value = 0.1132267775161799
$ws.Range("B3").Value = -0.1094296970220441
$ws.Range("C3").Value = -0.5919154831931398
$ws.Range("D3").Value = 0.346977032977319
$ws.Range("E3").Value = 0.6251071122536418
$ws.Range("G3").Value = 0.3523991642135738
$ws.Range("H3").Value = 0.08454266992975468
$ws.Range("I3").Value = 0.5651828635195373
$ws.Range("B4").Value = 0.5468644997283544
$ws.Range("C4").Value = 0.0006550715212962707
$ws.Range("D4").Value = 1.151798602178438
$ws.Range("B5").Value = 0.907150703699853
$ws.Range("C5").Value = 0.0001736751018040655
$ws.Range("D5").Value = 2.229846172697675

$ws = $wb.Worksheets.Item("Protists")
$ws.Range("B2").Value = -0.9275644248602527
$ws.Range("C2").Value = -2.021174175828682
$ws.Range("D2").Value = 0.1781684828495378
$ws.Range("E2").Value = 0.07497857754927173
$ws.Range("G2").Value = 0.0112702969106404
$ws.Range("H2").Value = 0.0000000001168747347216687
$ws.Range("I2").Value = 0.09208769738712898
$ws.Range("B3").Value = 0.01281899394636703
$ws.Range("C3").Value = -0.3222326932963732
$ws.Range("D3").Value = 0.3310391849863202
$ws.Range("E3").Value = 0.9284490145672666
$ws.Range("G3").Value = 0.3977535799041945
$ws.Range("H3").Value = 0.1175956927889487
$ws.Range("I3").Value = 0.6104745505746456
$ws.Range("B4").Value = 0.3772554572548459
$ws.Range("C4").Value = 0.00006384747985893152
$ws.Range("D4").Value = 0.8610929354152685
$ws.Range("B5").Value = 0.8183031899752338
$ws.Range("C5").Value = 0.001003705964852943
$ws.Range("D5").Value = 1.99096592208239

$ws = $wb.Worksheets.Item("Bacteria")
$ws.Range("B2").Value = -6.116724577228305
$ws.Range("C2").Value = -7.574624260268472
$ws.Range("D2").Value = -4.571076654326432
$ws.Range("E2").Value = 0.0004284490145671782
$ws.Range("G2").Value = 0.08847559383205764
$ws.Range("H2").Value = 0.000000131696401501414
$ws.Range("I2").Value = 0.2676651759927107
$ws.Range("B3").Value = -0.6386565865471112
$ws.Range("C3").Value = -1.497441672046309
$ws.Range("D3").Value = 0.1164906822615449
$ws.Range("E3").Value = 0.1135389888603255
$ws.Range("G3").Value = 0.2657690043821767
$ws.Range("H3").Value = 0.03286038096274363
$ws.Range("I3").Value = 0.4929916174201059
$ws.Range("B4").Value = 0.7532221810628891
$ws.Range("C4").Value = 0.0006079170835843788
$ws.Range("D4").Value = 1.697642035587546
$ws.Range("B5").Value = 0.857282876066121
$ws.Range("C5").Value = 0.000437916461978934
$ws.Range("D5").Value = 2.324973984042936
